# Update the NATMI ligand-receptor (Adm-Calcr) results sheet with the
# re-computed TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" sending-cluster pairs (old rows 10-11) are no longer
# part of the recomputed results, so those rows are removed entirely.
# Deleting row 10 twice removes both old row 10 and the row that shifts
# into its place (old row 11); the shared-string table and remaining row
# references are renumbered automatically by the engine.
$ws.Rows.Item(10).Delete() | Out-Null
$ws.Rows.Item(10).Delete() | Out-Null

# Refresh the recomputed numeric columns (G-T) for the remaining rows
# (now rows 2-9) to reflect the new TPM-derived values.

$ws.Range("G2").Value = 7.312510333333333
$ws.Range("H2").Value = 21.937531
$ws.Range("I2").Value = 0.2480018119509629
$ws.Range("J2").Value = 0.251002285750873
$ws.Range("M2").Value = 0.1045313333333333
$ws.Range("N2").Value = 0.313594
$ws.Range("O2").Value = 0.1482480886574569
$ws.Range("P2").Value = 0.2070264921079542
$ws.Range("Q2").Value = 0.764386455157111
$ws.Range("R2").Value = 6.879478096413999
$ws.Range("S2").Value = 0.0367657946053163
$ws.Range("T2").Value = 0.05196412273008157
$ws.Range("G3").Value = 7.312510333333333
$ws.Range("H3").Value = 21.937531
$ws.Range("I3").Value = 0.2480018119509629
$ws.Range("J3").Value = 0.251002285750873
$ws.Range("M3").Value = 0.6005795
$ws.Range("N3").Value = 1.201159
$ws.Range("O3").Value = 0.8517519113425431
$ws.Range("P3").Value = 0.7929735078920458
$ws.Range("Q3").Value = 4.391743799738167
$ws.Range("R3").Value = 26.350462798429
$ws.Range("S3").Value = 0.2112360173456466
$ws.Range("T3").Value = 0.1990381630207914
$ws.Range("I4").Value = 0.7058586226052527
$ws.Range("J4").Value = 0.7143985211120685
$ws.Range("M4").Value = 0.1045313333333333
$ws.Range("N4").Value = 0.313594
$ws.Range("O4").Value = 0.1482480886574569
$ws.Range("P4").Value = 0.2070264921079542
$ws.Range("Q4").Value = 2.175583985176667
$ws.Range("R4").Value = 19.58025586659
$ws.Range("S4").Value = 0.1046421916636139
$ws.Range("T4").Value = 0.1478994197929418
$ws.Range("I5").Value = 0.7058586226052527
$ws.Range("J5").Value = 0.7143985211120685
$ws.Range("M5").Value = 0.6005795
$ws.Range("N5").Value = 1.201159
$ws.Range("O5").Value = 0.8517519113425431
$ws.Range("P5").Value = 0.7929735078920458
$ws.Range("Q5").Value = 12.4997079857275
$ws.Range("R5").Value = 74.99824791436501
$ws.Range("S5").Value = 0.6012164309416388
$ws.Range("T5").Value = 0.5664991013191266
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.303044
$ws.Range("H6").Value = 0.909132
$ws.Range("I6").Value = 0.0102776553707253
$ws.Range("J6").Value = 0.01040200057377754
$ws.Range("M6").Value = 0.1045313333333333
$ws.Range("N6").Value = 0.313594
$ws.Range("O6").Value = 0.1482480886574569
$ws.Range("P6").Value = 0.2070264921079542
$ws.Range("Q6").Value = 0.03167759337866666
$ws.Range("R6").Value = 0.285098340408
$ws.Range("S6").Value = 0.001523642764590073
$ws.Range("T6").Value = 0.002153489689694092
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.303044
$ws.Range("H7").Value = 0.909132
$ws.Range("I7").Value = 0.0102776553707253
$ws.Range("J7").Value = 0.01040200057377754
$ws.Range("M7").Value = 0.6005795
$ws.Range("N7").Value = 1.201159
$ws.Range("O7").Value = 0.8517519113425431
$ws.Range("P7").Value = 0.7929735078920458
$ws.Range("Q7").Value = 0.182002013998
$ws.Range("R7").Value = 1.092012083988
$ws.Range("S7").Value = 0.008754012606135233
$ws.Range("T7").Value = 0.008248510884083451
$ws.Range("G8").Value = 1.057414
$ws.Range("H8").Value = 2.114828
$ws.Range("I8").Value = 0.03586191007305911
$ws.Range("J8").Value = 0.02419719256328104
$ws.Range("M8").Value = 0.1045313333333333
$ws.Range("N8").Value = 0.313594
$ws.Range("O8").Value = 0.1482480886574569
$ws.Range("P8").Value = 0.2070264921079542
$ws.Range("Q8").Value = 0.1105328953053333
$ws.Range("R8").Value = 0.663197371832
$ws.Range("S8").Value = 0.005316459623936613
$ws.Range("T8").Value = 0.005009459895236749
$ws.Range("G9").Value = 1.057414
$ws.Range("H9").Value = 2.114828
$ws.Range("I9").Value = 0.03586191007305911
$ws.Range("J9").Value = 0.02419719256328104
$ws.Range("M9").Value = 0.6005795
$ws.Range("N9").Value = 1.201159
$ws.Range("O9").Value = 0.8517519113425431
$ws.Range("P9").Value = 0.7929735078920458
$ws.Range("Q9").Value = 0.6350611714130001
$ws.Range("R9").Value = 2.540244685652
$ws.Range("S9").Value = 0.0305454504491225
$ws.Range("T9").Value = 0.01918773266804429
